$d = $word.ActiveDocument

# --- 1. Intro paragraph: "Clone the "Content Tools" Git repo" -> "Clone the GitHub repo" ---
# First, change the hyperlink display text.
$d.Content.Find.Execute("Content Tools", $true, $false, $false, $false, $false,
                         $true, 1, $false, "GitHub repo", 2)

# Now collapse the surrounding quote marks and the trailing "Git " word around the
# (now renamed) hyperlink text, leaving: Clone the GitHub repo to get the GTU source code.
$d.Content.Find.Execute([char]0x201C + "GitHub repo" + [char]0x201D + " Git repo to get",
                         $true, $false, $false, $false, $false,
                         $true, 1, $false, "GitHub repo to get", 2)

# --- 2. Background paragraph edits ---

# a) "GenTopics enumerates" -> "The original GenTopics application enumerates"
$d.Content.Find.Execute("file. GenTopics enumerates", $true, $false, $false, $false, $false,
                         $true, 1, $false, "file. The original GenTopics application enumerates", 2)

# b) "...stored in the OSG CPub source depot (SD) repository, and writers author content into them by using XMetaL." ->
#    "...stored in a source repository, and writers author content into them by using XMetaL or another XML editor."
$d.Content.Find.Execute("stored in the OSG CPub source depot (SD) repository, and writers author content into them by using XMetaL.",
                         $true, $false, $false, $false, $false,
                         $true, 1, $false,
                         "stored in a source repository, and writers author content into them by using XMetaL or another XML editor.", 2)

# c) "Writers are required to keep the SD topics in sync with source code changes manually, which is time-consuming and error-prone."
#    -> "Writers are required to keep the source topics in sync with source code changes manually, which can be time-consuming and error-prone."
$d.Content.Find.Execute("Writers are required to keep the SD topics in sync with source code changes manually, which is time-consuming and error-prone.",
                         $true, $false, $false, $false, $false,
                         $true, 1, $false,
                         "Writers are required to keep the source topics in sync with source code changes manually, which can be time-consuming and error-prone.", 2)

# --- 3. Move the _GoBack bookmark from the end of the document to just after
#        "...or another XML editor" (the last edit location), matching Word's
#        automatic tracking of the most recent edit point. ---
$old = $d.Bookmarks("_GoBack")
$old.Delete()

$d.Content.Find.Execute("or another XML editor", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 0)
$found = $d.Content.Find.Found
$r = $d.Content
$r.Find.Execute("or another XML editor", $true, $false, $false, $false, $false,
                $true, 1, $false, "", 0)
$anchor = $r.Duplicate
$anchor.Collapse(0)
$d.Bookmarks.Add("_GoBack", $anchor)
